# Apply corrected data values to the "données03" sheet.
# (Some rows had logic problems in column A / C that needed fixing.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("données03")

# Row 19
$ws.Range("A19").Value = 8.09
$ws.Range("C19").Value = 142

# Row 25
$ws.Range("A25").Value = 9.34
$ws.Range("C25").Value = 126

# Row 27
$ws.Range("A27").Value = 6.2700000000000005
$ws.Range("C27").Value = 128

# Row 32
$ws.Range("A32").Value = 11.01
$ws.Range("C32").Value = 132

# Row 33
$ws.Range("A33").Value = 15.39
$ws.Range("C33").Value = 122

# Row 38
$ws.Range("A38").Value = 31.319999999999997
$ws.Range("C38").Value = 135

# Row 39
$ws.Range("A39").Value = 64.239999999999995
$ws.Range("C39").Value = 139

# Row 40
$ws.Range("A40").Value = 14.219999999999999
$ws.Range("C40").Value = 134
